$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-122 down to 21-123
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new data record
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = (Get-Date -Year 2021 -Month 12 -Day 17 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112031
$ws.Range("G20").Value = "Poroto verde"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 15000
$ws.Range("N20").Value = '$/saco 25 kilos'
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 600
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
